# The GSC "Coverage" export rolled forward by one day: the oldest date
# (2025-10-15, row 2 of the "Chart" sheet) drops off the front of the
# rolling window and every later row shifts up to take its place.
# Deleting the entire row 2 reproduces exactly that: remaining rows
# (old row 3..89) shift up to become new rows 2..88, and the sheet's
# used range shrinks from A1:D89 to A1:D88 - matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows("2:2").Delete()
